$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency Price (col D) / Volume(1h) (col E) data scraped from
# coinranking.com. Rows 50-51 (Aave / FraxShare) also swap rank position,
# so their Coin (B) and Link (C) values are updated too.
#
# Numeric-looking Price values are written with a leading single-quote so
# Excel stores them as text (matching the source data), instead of parsing
# them into Number cells.

$ws.Cells.Item(2, 4).Value = '46.840.08'
$ws.Cells.Item(2, 5).Value = '  -0.22%  '
$ws.Cells.Item(3, 4).Value = '2.265.79'
$ws.Cells.Item(3, 5).Value = '  -3.61%  '
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 4).Value = '''298.62'
$ws.Cells.Item(5, 5).Value = '  -2.81%  '
$ws.Cells.Item(6, 4).Value = '''98.47'
$ws.Cells.Item(6, 5).Value = '  +0.20%  '
$ws.Cells.Item(7, 4).Value = '''0.574'
$ws.Cells.Item(7, 5).Value = '  -0.76%  '
$ws.Cells.Item(8, 5).Value = '  +0.16%  '
$ws.Cells.Item(9, 4).Value = '''0.504'
$ws.Cells.Item(9, 5).Value = '  -6.80%  '
$ws.Cells.Item(10, 4).Value = '''34.88'
$ws.Cells.Item(10, 5).Value = '  -2.95%  '
$ws.Cells.Item(11, 4).Value = '''0.0794'
$ws.Cells.Item(11, 5).Value = '  -1.69%  '
$ws.Cells.Item(12, 4).Value = '''7.01'
$ws.Cells.Item(12, 5).Value = '  -5.86%  '
$ws.Cells.Item(13, 5).Value = '  -1.68%  '
$ws.Cells.Item(14, 4).Value = '2.614.31'
$ws.Cells.Item(14, 5).Value = '  -3.47%  '
$ws.Cells.Item(15, 4).Value = '2.270.75'
$ws.Cells.Item(15, 5).Value = '  -3.46%  '
$ws.Cells.Item(16, 4).Value = '''13.62'
$ws.Cells.Item(16, 5).Value = '  -4.67%  '
$ws.Cells.Item(17, 4).Value = '46.838.20'
$ws.Cells.Item(17, 5).Value = '  +0.07%  '
$ws.Cells.Item(18, 4).Value = '''0.793'
$ws.Cells.Item(18, 5).Value = '  -4.68%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0970'
$ws.Cells.Item(19, 5).Value = '  +1.84%  '
$ws.Cells.Item(20, 4).Value = '''12.40'
$ws.Cells.Item(20, 5).Value = '  -9.23%  '
$ws.Cells.Item(21, 4).Value = '''5.79'
$ws.Cells.Item(21, 5).Value = '  -6.62%  '
$ws.Cells.Item(22, 4).Value = '''65.76'
$ws.Cells.Item(22, 5).Value = '  -1.86%  '
$ws.Cells.Item(23, 4).Value = '''245.12'
$ws.Cells.Item(23, 5).Value = '  -0.05%  '
$ws.Cells.Item(24, 4).Value = '''2.78'
$ws.Cells.Item(24, 5).Value = '  -7.40%  '
$ws.Cells.Item(25, 5).Value = '  +0.60%  '
$ws.Cells.Item(26, 4).Value = '''1.85'
$ws.Cells.Item(26, 5).Value = '  -7.86%  '
$ws.Cells.Item(27, 4).Value = '''41.32'
$ws.Cells.Item(27, 5).Value = '  -2.02%  '
$ws.Cells.Item(28, 4).Value = '''2.21'
$ws.Cells.Item(28, 5).Value = '  -3.46%  '
$ws.Cells.Item(29, 5).Value = '  -3.79%  '
$ws.Cells.Item(30, 4).Value = '''20.04'
$ws.Cells.Item(30, 5).Value = '  -0.96%  '
$ws.Cells.Item(31, 4).Value = '''2.82'
$ws.Cells.Item(31, 5).Value = '  +7.12%  '
$ws.Cells.Item(32, 4).Value = '''3.32'
$ws.Cells.Item(32, 5).Value = '  +4.85%  '
$ws.Cells.Item(33, 4).Value = '''145.08'
$ws.Cells.Item(33, 5).Value = '  -4.85%  '
$ws.Cells.Item(34, 4).Value = '''5.29'
$ws.Cells.Item(34, 5).Value = '  -8.04%  '
$ws.Cells.Item(35, 4).Value = '''0.0765'
$ws.Cells.Item(35, 5).Value = '  -6.17%  '
$ws.Cells.Item(36, 5).Value = '  +1.31%  '
$ws.Cells.Item(37, 5).Value = '  -2.80%  '
$ws.Cells.Item(38, 4).Value = '''15.39'
$ws.Cells.Item(38, 5).Value = '  +10.96%  '
$ws.Cells.Item(39, 4).Value = '''1.65'
$ws.Cells.Item(39, 5).Value = '  -9.94%  '
$ws.Cells.Item(40, 4).Value = '''3.81'
$ws.Cells.Item(40, 5).Value = '  -6.57%  '
$ws.Cells.Item(41, 4).Value = '''0.0293'
$ws.Cells.Item(41, 5).Value = '  -7.22%  '
$ws.Cells.Item(42, 5).Value = '  -11.46%  '
$ws.Cells.Item(43, 5).Value = '  +0.06%  '
$ws.Cells.Item(44, 4).Value = '''93.73'
$ws.Cells.Item(44, 5).Value = '  +15.90%  '
$ws.Cells.Item(45, 4).Value = '1.782.24'
$ws.Cells.Item(45, 5).Value = '  -7.87%  '
$ws.Cells.Item(46, 4).Value = '''1.87'
$ws.Cells.Item(46, 5).Value = '  -5.56%  '
$ws.Cells.Item(47, 4).Value = '''70.39'
$ws.Cells.Item(47, 5).Value = '  -5.45%  '
$ws.Cells.Item(48, 5).Value = '  -7.68%  '
$ws.Cells.Item(49, 5).Value = '  -3.39%  '
$ws.Cells.Item(50, 2).Value = 'FraxShare'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(50, 4).Value = '''7.88'
$ws.Cells.Item(50, 5).Value = '  -1.92%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '''93.91'
$ws.Cells.Item(51, 5).Value = '  -5.35%  '
